$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.373.40'
$ws.Range('E2').Value = '  +2.63%  '

# Row 3
$ws.Range('D3').Value = '2.109.57'
$ws.Range('E3').Value = '  +1.05%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.37%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.75%  '

# Row 6
$ws.Range('E6').Value = '  -0.33%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5234'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.75%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4444'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.55%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.54'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.41%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09464'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.36%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.176'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.31%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.96%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.750'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.79%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.932'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.00%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.096.67'
$ws.Range('E15').Value = '  -1.42%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.56%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001167'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.68%  '

# Row 18
$ws.Range('E18').Value = '  -0.35%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.67%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06722'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.24%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.324'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.64%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '

# Row 23
$ws.Range('D23').Value = '30.410.45'
$ws.Range('E23').Value = '  +2.55%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.67%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.314'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.15%  '

# Row 26
$ws.Range('D26').Value = '2.378.00'
$ws.Range('E26').Value = '  +1.49%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.28%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.548'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.31%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.58%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.84%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.150'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.54%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.761'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.89%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1056'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.15%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.891'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.38%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.277'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.41%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.925'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.13%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.52'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.17%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02637'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.03%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06813'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.01%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.7045'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.12%  '

# Row 41
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.77%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.347'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.11%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2233'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.21%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6850'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.87%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.41%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.366'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.46%  '

# Row 47
$ws.Range('E47').Value = '  -0.25%  '

# Row 48
$ws.Range('E48').Value = '  +16.04%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.647'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.09%  '

# Row 50
$ws.Range('E50').Value = '  +2.25%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.225'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.89%  '
